$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "Hello world"
$ws.Range("F4").Value = "Hello, world"
